$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.880.75'
$ws.Range("D3").Value = '1.730.28'
$ws.Range("D4").Value = "'0.9972"
$ws.Range("D5").Value = "'241.98"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = "'0.9977"
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("D7").Value = "'0.4916"
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("D8").Value = "'0.2601"
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").Value = "'0.06220"
$ws.Range("E9").Value = '  +0.32%  '
$ws.Range("D10").Value = '1.736.65'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = "'16.06"
$ws.Range("E11").Value = '  +3.49%  '
$ws.Range("D12").Value = "'0.06904"
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").Value = "'0.6105"
$ws.Range("E13").Value = '  +1.86%  '
$ws.Range("D14").Value = "'4.502"
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D15").Value = "'77.36"
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("D17").Value = '26.868.34'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").Value = "'0.9972"
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").Value = "'0.000007188"
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = "'11.46"
$ws.Range("E20").Value = '  +0.99%  '
$ws.Range("D21").Value = '1.959.05'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("D22").Value = "'4.437"
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = "'8.562"
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("D24").Value = "'5.124"
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("D25").Value = "'138.81"
$ws.Range("E25").Value = '  +0.90%  '
$ws.Range("D26").Value = "'15.34"
$ws.Range("E26").Value = '  +0.76%  '
$ws.Range("D27").Value = "'1.789"
$ws.Range("E27").Value = '  +5.26%  '
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").Value = "'106.32"
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -0.19%  '
$ws.Range("D31").Value = "'0.07995"
$ws.Range("E31").Value = '  +0.57%  '
$ws.Range("D32").Value = "'3.681"
$ws.Range("E32").Value = '  +0.01%  '
$ws.Range("D33").Value = "'0.04533"
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").Value = "'0.9969"
$ws.Range("D35").Value = "'2.610"
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").Value = "'1.010"
$ws.Range("E36").Value = '  +1.39%  '
$ws.Range("D37").Value = "'0.6248"
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").Value = "'0.9332"
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").Value = "'2.058"
$ws.Range("E39").Value = '  +5.08%  '
$ws.Range("D40").Value = "'2.440"
$ws.Range("E40").Value = '  +1.95%  '
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("D43").Value = "'5.673"
$ws.Range("E43").Value = '  +4.30%  '
$ws.Range("D44").Value = "'99.76"
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("D45").Value = "'0.3869"
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("D46").Value = "'6.957"
$ws.Range("E46").Value = '  +3.66%  '
$ws.Range("D48").Value = "'0.05389"
$ws.Range("D49").Value = "'7.942"
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("D50").Value = "'30.26"
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = "'1.240"
$ws.Range("E51").Value = '  +0.01%  '
